$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Udemy")

# ------------------------------------------------------------------
# The sheet has a list of course sections, each one a vertically
# merged cell in column A (the section title) next to a column B
# list of lessons, bordered top/middle/bottom via styles 13/14/15.
# We are inserting a new section - "Node JS + Express + TS" - right
# where the old, bare "Node JS + Express + TS" / "Course Roundup"
# title rows used to sit (rows 210/211), pushing "Course Roundup"
# down to row 220 and growing the sheet from A1:B211 to A1:B220.
# ------------------------------------------------------------------

# Preserve the "Course Roundup" row (currently A211) and relocate it
# to its new home at A220, keeping both its text and its formatting.
$roundupText = $ws.Range("A211").Value2
$ws.Range("A211").Copy()
$ws.Range("A220").PasteSpecial(-4122)
$ws.Range("A220").Value2 = $roundupText

# Stamp the new section's formatting by copying it from the
# "ReactJS and TS" block (A196:B208), which has the identical
# top/middle/bottom row styling we need, just re-tiled to the 9 rows
# (210-218) this section needs instead of 13.
$ws.Range("A196:B196").Copy()
$ws.Range("A210:B210").PasteSpecial(-4122)

$ws.Range("A197:B197").Copy()
$ws.Range("A211:B217").PasteSpecial(-4122)

$ws.Range("A208:B208").Copy()
$ws.Range("A218:B218").PasteSpecial(-4122)

# The old A211 cell content has already been carried over to A220,
# clear it out before reusing as part of the new section.
$ws.Range("A211").ClearContents()

# Fill in the new section's text.
$ws.Range("A210").Value2 = "Node JS + Express + TS"
$ws.Range("B210").Value2 = "Module Introduction"
$ws.Range("B211").Value2 = "Executing Typescript Code with Node.JS"
$ws.Range("B212").Value2 = "Setting Up a Project"
$ws.Range("B213").Value2 = "Finished Setup & Working with Types (In NodeJS + Express)"
$ws.Range("B214").Value2 = "Adding Middleware and Types"
$ws.Range("B215").Value2 = "Working with Controllers and Parsing with Request Bodies"
$ws.Range("B216").Value2 = "More CRUD Operations"
$ws.Range("B217").Value2 = "Wrap Up"
$ws.Range("B218").Value2 = "Useful Resources and Links"

# Merge the section title cell, just like every other section.
$ws.Range("A210:A218").Merge()

# Match the workbook's saved selection/view state to the new layout.
$ws.Range("A210:B218").Select()
